# Update Name of Algo
# Applies the updated numeric values produced by re-running the algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5.602600000000001
$ws.Range("A12").Value = -21.38819999999999
$ws.Range("C14").Value = -13.7001
$ws.Range("C19").Value = -12.90090000000001
$ws.Range("B23").Value = 8.410000000000002
$ws.Range("C24").Value = -13.06759999999999
$ws.Range("A27").Value = -21.98289999999999
$ws.Range("B28").Value = 5.305099999999998
$ws.Range("A32").Value = -21.04469999999998
$ws.Range("B32").Value = 6.160599999999997
$ws.Range("B34").Value = 9.497400000000004
$ws.Range("A36").Value = -19.7297
$ws.Range("A38").Value = -19.8748
$ws.Range("C38").Value = -10.68620000000001
$ws.Range("C41").Value = -12.66140000000001
$ws.Range("B42").Value = 10.0325
$ws.Range("A46").Value = -21.72689999999999
$ws.Range("B49").Value = 5.2874
$ws.Range("C52").Value = -11.195
$ws.Range("A54").Value = -21.92180000000001
$ws.Range("B54").Value = 4.2381
$ws.Range("A55").Value = -22.0343
$ws.Range("A56").Value = -22.16230000000001
$ws.Range("A67").Value = -21.52949999999997
$ws.Range("A69").Value = -21.60469999999997
$ws.Range("A72").Value = -22.05099999999999
$ws.Range("C72").Value = -12.2692
$ws.Range("B78").Value = 9.6401
$ws.Range("C78").Value = -12.3019
$ws.Range("B80").Value = 9.428600000000003
$ws.Range("A83").Value = -21.60439999999999
$ws.Range("C83").Value = -13.00339999999999
$ws.Range("C85").Value = -12.98079999999999
$ws.Range("A86").Value = -21.45249999999999
$ws.Range("C86").Value = -12.828
$ws.Range("C90").Value = -10.0916
$ws.Range("A91").Value = -20.70319999999999
$ws.Range("A93").Value = -21.35350000000001
$ws.Range("C96").Value = -10.1033
$ws.Range("B97").Value = 6.298399999999995
$ws.Range("A99").Value = -21.97490000000001
$ws.Range("B99").Value = 5.166499999999997
$ws.Range("B101").Value = 4.026800000000002
$ws.Range("C103").Value = -12.44849999999999
$ws.Range("A104").Value = -21.4967
